# PBL.xlsx update:
# - PBI [2] (CSV import) marked as Done
# - New PBI [9] "Sprint-Daten editieren, speichern und wieder laden" added
# - PBI [5] "Sprint-Daten CSV import" moved up (next PBI selected)
# - Table re-ordered accordingly; PBIs [3] and [4] moved to the bottom

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C texts (wrapped Akzeptanzkriterien / notes)
$cEclipseGit    = "Akzeptanzkriterien:`n- Eclipse-Projekt`n- Git-Repo"
$cIdTitleSum    = "Akzeptanzkriterien:`n- Id, Title, Summary und Estimate"
$cSprintFields  = "Akzeptanzkriterien:`n- Sprintname, Anfangsdatum, Enddatum, geplante Kapazität, geplanter Aufwand, tatsächliche Kapazität, Aufwand Done"
$cStatusCsv     = "Akzeptanzkriterien:`n- aus dem CSV laden, im PBL anzeigen, im Burnup anzeigen und persistieren.`n- Status: Todo, In Progress, Done, Canceled"
$cTrendKapazit  = "Akzeptanzkriterien:`n- Trendlinie gerücksichtigt die Kapazität."
$cSprintEdit    = "Akzeptanzkriterien:`n- Alle Felder aus [5] können editiert werden.`n- Sprints können gelöscht werden.`n- Neue Sprints können hinzugefügt werden.`n- Die Persistierung erfolgt als XML-Datei."
$cZeitstempel   = "Akzeptanzkriterien:`n- Als Zeitstempel wird das Datum der CSV-Datei verwendet."
$cOffenePunkte  = "Offene Punkte:`n- Wohin persistieren?"

# Row definitions for the new table layout (rows 2..10)
$rows = @(
    @{ Row=2;  A=1; B="Als Entwickler möchte ich einen initialen Projektsetup haben"; C=$cEclipseGit;   D=0.5; E="Done"; H=45 },
    @{ Row=3;  A=2; B="Als PO möchte ich PBIs aus einer CSV-Datei importieren und als PBL in einer Tabelle anzeigen können"; C=$cIdTitleSum; D=3; E="Done"; H=45 },
    @{ Row=4;  A=5; B="Als PO möchte ich Sprint-Daten aus einer CSV-Datei importieren und tabellarisch anzeigen können."; C=$cSprintFields; D=2; E="Todo"; H=90 },
    @{ Row=5;  A=6; B="Als PO möchte ich den Status eines PBIs verwenden können"; C=$cStatusCsv; D=2; E="Todo"; H=90 },
    @{ Row=6;  A=7; B="Als PO möchte ich eine Trendlinie mit der durchschnittlichen Geschwindigkeit im Burnup anzeigen lassen können"; C=$cTrendKapazit; D=2; E="Todo"; H=45 },
    @{ Row=7;  A=8; B="Als PO möchte ich Trendlinien mit minimaler und maximaler Geschwindigkeit im Burnup anzeigen lassen können"; C=$null; D=1; E="Todo"; H=15 },
    @{ Row=8;  A=9; B="Als PO möchte ich Sprint-Daten editieren, speichern und wieder laden können."; C=$cSprintEdit; D=3; E="Todo"; H=135 },
    @{ Row=9;  A=3; B="Als PO möchte ich mehrere Versionen des PBLs als CSV laden und den zeitlichen Verlauf des Gesamtaufwands darstellen können."; C=$cZeitstempel; D=2; E="Todo"; H=60 },
    @{ Row=10; A=4; B="Als PO möchte ich mehrere zeitliche Versionen des PBLs persistieren und wieder laden können"; C=$cOffenePunkte; D=3; E="Todo"; H=30 }
)

foreach ($r in $rows) {
    $ri = $r.Row
    $ws.Cells.Item($ri, 1).Value2 = $r.A
    $ws.Cells.Item($ri, 2).Value2 = $r.B
    $cCell = $ws.Cells.Item($ri, 3)
    if ($r.C -ne $null) {
        $cCell.Value2 = $r.C
        $cCell.WrapText = $true
    } else {
        $cCell.Clear() | Out-Null
    }
    $ws.Cells.Item($ri, 4).Value2 = $r.D
    $ws.Cells.Item($ri, 5).Value2 = $r.E
    $ws.Rows.Item($ri).RowHeight = $r.H
}

# Column A width adjustment (bestFit 2.7109375 -> 3)
$ws.Columns.Item(1).ColumnWidth = 2.1666667

# Update dimension / selection to match new extent
$ws.Range("D9").Select() | Out-Null
